$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Locate the FIRST "Change the values so that the total data rate..."
#    paragraph (there are two identical ones; only the first is restructured
#    into three runs - with the middle one underlined - and gets a new
#    "Note:" paragraph appended right after it).
# ---------------------------------------------------------------------------
$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("Change the values so that the total data rate offered to the channel (that is sum")) {
        $targetIdx = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIdx)
$r = $p.Range

$part1 = "Change the values so that the total data rate offered to the channel ("
$part2 = "that is sum of data rates of all CBR flows"
$part3 = [char]0x29 + " is about 10% of the channel data rate, equally divided among all sources (keep " + [char]0x201C + "dataRate" + [char]0x201D + " for all sources the same).  Keep increasing offered load to 20%, 30% " + [char]0x2026 + " 90%. What" + [char]0x2019 + "s the total throughput achieved for each value? What is the trend of the throughput vs offered load? Plot the values. "

$r.Text = $part1 + $part2 + $part3

$pStart = $r.Start
$s1 = $pStart
$e1 = $s1 + $part1.Length
$s2 = $e1
$e2 = $s2 + $part2.Length

$underlineRange = $d.Range($s2, $e2)
$underlineRange.Font.Underline = 1

# ---------------------------------------------------------------------------
# 2) Insert the new "Note:" paragraph right after the paragraph above.
# ---------------------------------------------------------------------------
$r.InsertParagraphAfter()
$notePara = $d.Paragraphs.Item($targetIdx + 1)
$notePara.Style = "Normal"
$notePara.Format.FirstLineIndent = 0
$notePara.Format.LeftIndent = 72

$notePart1 = "Note:"
$notePart2 = " You should not touch the " + [char]0x201C + "payload" + [char]0x201D + " value. Instead change offered load using the " + [char]0x201C + "dataRate" + [char]0x201D + " parameter."

$noteRange = $notePara.Range
$noteRange.Text = $notePart1 + $notePart2

$noteStart = $noteRange.Start
$boldRange = $d.Range($noteStart, $noteStart + $notePart1.Length)
$boldRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# 3) "Find out the maximum throughput possible..." -> "...total throughput..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Find out the maximum throughput possible when there is no contention (only one source). ", $true, $false, $false, $false, $false, $true, 1, $false, "Find out the total throughput possible when there is no contention (only one source). ", 2)

# ---------------------------------------------------------------------------
# 4) Remaining three "...What's the maximum throughput achieved..." ->
#    "...What's the total throughput achieved..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("maximum throughput achieved", $true, $false, $false, $false, $false, $true, 1, $false, "total throughput achieved", 2)

# ---------------------------------------------------------------------------
# 5) Add a (default) footer to the document's single section.
# ---------------------------------------------------------------------------
$wdHeaderFooterPrimary = 1
$footer = $d.Sections.Item(1).Footers.Item($wdHeaderFooterPrimary)
$footer.Range.Text = ""
